# AVIV_CheckoutPage test-data sheet refresh: new random Address / PhNo / OrderNo values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVIV_CheckoutPage")

# D6 / D8 hold long digit-only strings (Address, PhNo test data) that must stay TEXT
# (shared-string) cells, not be auto-converted to numbers. We stage the value with a
# leading apostrophe in a scratch cell, then paste-special *values only* onto the
# target cell so the text type carries over without disturbing the cell's existing
# (default) style the way setting NumberFormat directly would.
$helper = $ws.Range("ZZ1")

$helper.Value = "'6076706315"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$helper.Clear()

$helper.Value = "'8481652086"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$helper.Clear()

# D13 (OrderNo) is plain text already, no numeric coercion risk.
$ws.Range("D13").Value = "ORDER NUMBER: 1033"
